$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Harvard case classification: rename header columns and shift
# average_doctor / average_doctor_old data, refreshing underlying stats.

# --- Header swap (row 1) ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"


# --- Row 4 ---
$ws.Range("AI4").Value = 0.401
$ws.Range("AJ4").Value = 0.093
$ws.Range("AK4").Value = 0.306
$ws.Range("AU4").Value = 0.239
$ws.Range("AW4").Value = 0.161
$ws.Range("BA4").Value = 2.064
$ws.Range("BB4").Value = 0.143
$ws.Range("BC4").Value = 0.379
$ws.Range("BG4").Value = 0.729
$ws.Range("BH4").Value = 0.141
$ws.Range("BI4").Value = 0.375
$ws.Range("BM4").Value = 0.756
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.6879999999999999
$ws.Range("BQ4").Value = 0.765
$ws.Range("E4").Value = 0.478
$ws.Range("F4").Value = 0.052
$ws.Range("G4").Value = 0.229
$ws.Range("N4").Value = 0.474
$ws.Range("O4").Value = 0.061
$ws.Range("P4").Value = 0.248
$ws.Range("Q4").Value = 0.051
$ws.Range("R4").Value = 0.034
$ws.Range("S4").Value = 0.186
$ws.Range("W4").Value = 0.366

# --- Row 5 ---
$ws.Range("AI5").Value = 0.402
$ws.Range("AJ5").Value = 0.09
$ws.Range("AK5").Value = 0.3
$ws.Range("AU5").Value = 0.445
$ws.Range("AV5").Value = 0.077
$ws.Range("AW5").Value = 0.278
$ws.Range("BA5").Value = 1.3
$ws.Range("BB5").Value = 0.07199999999999999
$ws.Range("BC5").Value = 0.268
$ws.Range("BG5").Value = 0.382
$ws.Range("BH5").Value = 0.05
$ws.Range("BI5").Value = 0.223
$ws.Range("BM5").Value = 0.521
$ws.Range("BN5").Value = 0.047
$ws.Range("BO5").Value = 0.216
$ws.Range("BP5").Value = 0.433
$ws.Range("BQ5").Value = 0.456
$ws.Range("E5").Value = 0.599
$ws.Range("F5").Value = 0.059
$ws.Range("G5").Value = 0.242
$ws.Range("N5").Value = 0.716
$ws.Range("O5").Value = 0.076
$ws.Range("P5").Value = 0.275
$ws.Range("Q5").Value = 0.034
$ws.Range("R5").Value = 0.014
$ws.Range("S5").Value = 0.119
$ws.Range("W5").Value = 0.335
$ws.Range("X5").Value = 0.097
$ws.Range("Y5").Value = 0.312

# --- Row 6 ---
$ws.Range("AI6").Value = 0.401
$ws.Range("AU6").Value = 0.311
$ws.Range("BA6").Value = 1.589
$ws.Range("BG6").Value = 0.501
$ws.Range("BM6").Value = 0.617
$ws.Range("BP6").Value = 0.53
$ws.Range("BQ6").Value = 0.569
$ws.Range("E6").Value = 0.532
$ws.Range("N6").Value = 0.57
$ws.Range("Q6").Value = 0.041
$ws.Range("W6").Value = 0.35

# --- Row 7 ---
$ws.Range("AI7").Value = 0.402
$ws.Range("AU7").Value = 0.38
$ws.Range("BA7").Value = 1.402
$ws.Range("BG7").Value = 0.422
$ws.Range("BM7").Value = 0.556
$ws.Range("BP7").Value = 0.467
$ws.Range("BQ7").Value = 0.495
$ws.Range("E7").Value = 0.57
$ws.Range("N7").Value = 0.65
$ws.Range("Q7").Value = 0.036
$ws.Range("W7").Value = 0.341

# --- Row 8 ---
$ws.Range("AI8").Value = 0.472
$ws.Range("AJ8").Value = 0.137
$ws.Range("AK8").Value = 0.37
$ws.Range("AU8").Value = 0.394
$ws.Range("AV8").Value = 0.083
$ws.Range("AW8").Value = 0.288
$ws.Range("BA8").Value = 1.771
$ws.Range("BB8").Value = 0.108
$ws.Range("BC8").Value = 0.328
$ws.Range("BG8").Value = 0.5639999999999999
$ws.Range("BH8").Value = 0.108
$ws.Range("BI8").Value = 0.329
$ws.Range("BM8").Value = 0.675
$ws.Range("BN8").Value = 0.061
$ws.Range("BO8").Value = 0.247
$ws.Range("BP8").Value = 0.59
$ws.Range("BQ8").Value = 0.625
$ws.Range("E8").Value = 0.706
$ws.Range("F8").Value = 0.07000000000000001
$ws.Range("G8").Value = 0.265
$ws.Range("N8").Value = 0.802
$ws.Range("O8").Value = 0.058
$ws.Range("P8").Value = 0.242
$ws.Range("Q8").Value = 0.037
$ws.Range("W8").Value = 0.405
$ws.Range("X8").Value = 0.119
$ws.Range("Y8").Value = 0.345

# --- Row 9 ---
$ws.Range("AI9").Value = 0.419
$ws.Range("AJ9").Value = 0.243
$ws.Range("AK9").Value = 0.493
$ws.Range("BA9").Value = 1.744
$ws.Range("BB9").Value = 0.25
$ws.Range("BC9").Value = 0.5
$ws.Range("BG9").Value = 0.605
$ws.Range("BH9").Value = 0.239
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.651
$ws.Range("BN9").Value = 0.227
$ws.Range("BO9").Value = 0.477
$ws.Range("BP9").Value = 0.581
$ws.Range("BQ9").Value = 0.619
$ws.Range("E9").Value = 0.674
$ws.Range("F9").Value = 0.22
$ws.Range("G9").Value = 0.469
$ws.Range("N9").Value = 0.744
$ws.Range("O9").Value = 0.19
$ws.Range("P9").Value = 0.436
$ws.Range("W9").Value = 0.302
$ws.Range("X9").Value = 0.211
$ws.Range("Y9").Value = 0.459

# --- Row 10 ---
$ws.Range("AI10").Value = 0.512
$ws.Range("AJ10").Value = 0.25
$ws.Range("AK10").Value = 0.5
$ws.Range("AU10").Value = 0.395
$ws.Range("AV10").Value = 0.239
$ws.Range("AW10").Value = 0.489
$ws.Range("BA10").Value = 2.186
$ws.Range("BB10").Value = 0.211
$ws.Range("BC10").Value = 0.459
$ws.Range("BG10").Value = 0.674
$ws.Range("BH10").Value = 0.22
$ws.Range("BI10").Value = 0.469
$ws.Range("BM10").Value = 0.8139999999999999
$ws.Range("BN10").Value = 0.151
$ws.Range("BO10").Value = 0.389
$ws.Range("BP10").Value = 0.729
$ws.Range("BQ10").Value = 0.762
$ws.Range("E10").Value = 0.8139999999999999
$ws.Range("F10").Value = 0.151
$ws.Range("G10").Value = 0.389
$ws.Range("N10").Value = 0.93
$ws.Range("O10").Value = 0.065
$ws.Range("P10").Value = 0.255
$ws.Range("W10").Value = 0.512
$ws.Range("X10").Value = 0.25
$ws.Range("Y10").Value = 0.5

# --- Row 11 ---
$ws.Range("AI11").Value = 0.581
$ws.Range("AJ11").Value = 0.243
$ws.Range("AK11").Value = 0.493
$ws.Range("AU11").Value = 0.5580000000000001
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.497
$ws.Range("BA11").Value = 2.186
$ws.Range("BB11").Value = 0.211
$ws.Range("BC11").Value = 0.459
$ws.Range("BG11").Value = 0.674
$ws.Range("BH11").Value = 0.22
$ws.Range("BI11").Value = 0.469
$ws.Range("BM11").Value = 0.8139999999999999
$ws.Range("BN11").Value = 0.151
$ws.Range("BO11").Value = 0.389
$ws.Range("BP11").Value = 0.729
$ws.Range("BQ11").Value = 0.768
$ws.Range("E11").Value = 0.86
$ws.Range("F11").Value = 0.12
$ws.Range("G11").Value = 0.347
$ws.Range("N11").Value = 0.93
$ws.Range("O11").Value = 0.065
$ws.Range("P11").Value = 0.255
$ws.Range("W11").Value = 0.512
$ws.Range("X11").Value = 0.25
$ws.Range("Y11").Value = 0.5

# --- Row 12 ---
$ws.Range("AI12").Value = 1.6
$ws.Range("AJ12").Value = 1.44
$ws.Range("AK12").Value = 1.2
$ws.Range("AU12").Value = 2.846
$ws.Range("AV12").Value = 3.361
$ws.Range("AW12").Value = 1.833
$ws.Range("BA12").Value = 3.767
$ws.Range("BB12").Value = 0.44
$ws.Range("BC12").Value = 0.663
$ws.Range("BG12").Value = 1.138
$ws.Range("BH12").Value = 0.188
$ws.Range("BI12").Value = 0.433
$ws.Range("BM12").Value = 1.229
$ws.Range("BN12").Value = 0.233
$ws.Range("BO12").Value = 0.483
$ws.Range("BP12").Value = 1.256
$ws.Range("BQ12").Value = 1.243
$ws.Range("E12").Value = 1.405
$ws.Range("F12").Value = 0.836
$ws.Range("G12").Value = 0.914
$ws.Range("N12").Value = 1.25
$ws.Range("O12").Value = 0.287
$ws.Range("P12").Value = 0.536
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.432
$ws.Range("Y12").Value = 0.657

# --- Row 13 ---
$ws.Range("AI13").Value = 1.154
$ws.Range("AJ13").Value = 0.303
$ws.Range("AK13").Value = 0.551
$ws.Range("AU13").Value = 2.039
$ws.Range("AV13").Value = 0.339
$ws.Range("AW13").Value = 0.582
$ws.Range("BA13").Value = 2.159
$ws.Range("BB13").Value = 0.277
$ws.Range("BC13").Value = 0.527
$ws.Range("BG13").Value = 0.542
$ws.Range("BH13").Value = 0.05
$ws.Range("BI13").Value = 0.224
$ws.Range("BM13").Value = 0.776
$ws.Range("BN13").Value = 0.164
$ws.Range("BO13").Value = 0.404
$ws.Range("BP13").Value = 0.72
$ws.Range("BQ13").Value = 0.661
$ws.Range("E13").Value = 1.405
$ws.Range("F13").Value = 0.292
$ws.Range("G13").Value = 0.54
$ws.Range("N13").Value = 1.737
$ws.Range("O13").Value = 0.466
$ws.Range("P13").Value = 0.6830000000000001
$ws.Range("W13").Value = 0.985
$ws.Range("X13").Value = 0.194
$ws.Range("Y13").Value = 0.441
